$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3606.0588
$ws.Range("I17").Value = 8563
$ws.Range("J17").Value = 2945.1333
$ws.Range("K17").Value = 25689
$ws.Range("L17").Value = 8835.3999
$ws.Range("M17").Value = -25521
$ws.Range("N17").Value = -9171.3999
$ws.Range("H42").Value = 230.4
$ws.Range("I42").Value = 160.75
$ws.Range("K42").Value = 482.25
$ws.Range("M42").Value = -252.25
$ws.Range("H74").Value = 3449.2222
$ws.Range("I74").Value = 2408.8
$ws.Range("K74").Value = 2408.8
$ws.Range("M74").Value = -1472.8
$ws.Range("H77").Value = 3449.2222
$ws.Range("I77").Value = 2408.8
$ws.Range("K77").Value = 12044
$ws.Range("M77").Value = -7364
$ws.Range("H106").Value = 2365.6155
$ws.Range("I106").Value = 2569
$ws.Range("K106").Value = 2569
$ws.Range("M106").Value = -1938
$ws.Range("H132").Value = 1601.5
$ws.Range("I132").Value = 1766.6666
$ws.Range("K132").Value = 5299.9998
$ws.Range("M132").Value = -2769.9998
$ws.Range("H137").Value = 1843.125
$ws.Range("I137").Value = 1124.1666
$ws.Range("K137").Value = 3372.4998
$ws.Range("M137").Value = -822.4998000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4529.451
$ws.Range("I32").Value = 3745.1633
$ws.Range("K32").Value = 3745.1633
$ws.Range("M32").Value = -3458.1633
$ws.Range("H37").Value = 20000
$ws.Range("J37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("N37").Value = -20546
$ws.Range("H74").Value = 1183.3636
$ws.Range("I74").Value = 457.6111
$ws.Range("J74").Value = 4449.25
$ws.Range("K74").Value = 457.6111
$ws.Range("L74").Value = 4449.25
$ws.Range("M74").Value = 416.3889
$ws.Range("N74").Value = -6197.25
$ws.Range("H77").Value = 1183.3636
$ws.Range("I77").Value = 457.6111
$ws.Range("J77").Value = 4449.25
$ws.Range("K77").Value = 2288.0555
$ws.Range("L77").Value = 22246.25
$ws.Range("M77").Value = 2079.9445
$ws.Range("N77").Value = -30982.25
$ws.Range("H109").Value = 58656.75
$ws.Range("J109").Value = 58656.75
$ws.Range("L109").Value = 58656.75
$ws.Range("N109").Value = -61430.75
$ws.Range("H132").Value = 1676.8334
$ws.Range("I132").Value = 1376.7858
$ws.Range("J132").Value = 2727
$ws.Range("K132").Value = 4130.357400000001
$ws.Range("L132").Value = 8181
$ws.Range("M132").Value = -1600.357400000001
$ws.Range("N132").Value = -13241

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 128647.31
$ws.Range("I86").Value = 4113.0835
$ws.Range("K86").Value = 4113.0835
$ws.Range("M86").Value = -2990.0835
$ws.Range("H89").Value = 128647.31
$ws.Range("I89").Value = 4113.0835
$ws.Range("K89").Value = 20565.4175
$ws.Range("M89").Value = -14949.4175
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("H107").Value = 950.1111
$ws.Range("I107").Value = 728.2222
$ws.Range("K107").Value = 728.2222
$ws.Range("M107").Value = 1191.7778
$ws.Range("H134").Value = 4437.514
$ws.Range("I134").Value = 4849.1787
$ws.Range("J134").Value = 2790.8572
$ws.Range("K134").Value = 14547.5361
$ws.Range("L134").Value = 8372.571599999999
$ws.Range("M134").Value = -12012.5361
$ws.Range("N134").Value = -13442.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2202.4707
$ws.Range("J31").Value = 2520
$ws.Range("L31").Value = 2520
$ws.Range("N31").Value = -3110
$ws.Range("H34").Value = 2202.4707
$ws.Range("J34").Value = 2520
$ws.Range("L34").Value = 2520
$ws.Range("N34").Value = -2924
$ws.Range("H53").Value = 48187.25
$ws.Range("J53").Value = 48187.25
$ws.Range("L53").Value = 48187.25
$ws.Range("N53").Value = -49401.25
$ws.Range("H58").Value = 7249576.5
$ws.Range("J58").Value = 4500
$ws.Range("L58").Value = 4500
$ws.Range("N58").Value = -4906
$ws.Range("H62").Value = 3729.8
$ws.Range("I62").Value = 3499.5
$ws.Range("J62").Value = 3883.3333
$ws.Range("K62").Value = 3499.5
$ws.Range("L62").Value = 3883.3333
$ws.Range("M62").Value = -2875.5
$ws.Range("N62").Value = -5131.3333
$ws.Range("H65").Value = 3729.8
$ws.Range("I65").Value = 3499.5
$ws.Range("J65").Value = 3883.3333
$ws.Range("K65").Value = 17497.5
$ws.Range("L65").Value = 19416.6665
$ws.Range("M65").Value = -14377.5
$ws.Range("N65").Value = -25656.6665
$ws.Range("H94").Value = 1417
$ws.Range("J94").Value = 1290.6666
$ws.Range("L94").Value = 1290.6666
$ws.Range("N94").Value = -2192.6666
$ws.Range("H122").Value = 4057.7273
$ws.Range("I122").Value = 2952.75
$ws.Range("K122").Value = 8858.25
$ws.Range("M122").Value = -6408.25
$ws.Range("H132").Value = 3134.5
$ws.Range("I132").Value = 1859.375
$ws.Range("K132").Value = 5578.125
$ws.Range("M132").Value = -3048.125
$ws.Range("H136").Value = 7249576.5
$ws.Range("J136").Value = 4500
$ws.Range("L136").Value = 13500
$ws.Range("N136").Value = -18600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 150
$ws.Range("K4").Value = 450
$ws.Range("M4").Value = -338
$ws.Range("H121").Value = 646
$ws.Range("I121").Value = 410
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 1230
$ws.Range("L121").Value = 3000
$ws.Range("M121").Value = 80
$ws.Range("N121").Value = -5620
$ws.Range("H124").Value = 6000
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H131").Value = 13424.825
$ws.Range("J131").Value = 14527.655
$ws.Range("L131").Value = 43582.965
$ws.Range("N131").Value = -53662.965

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 615.69696
$ws.Range("I97").Value = 614.8276
$ws.Range("K97").Value = 614.8276
$ws.Range("M97").Value = -118.8276
$ws.Range("H122").Value = 2169
$ws.Range("I122").Value = 2136.7
$ws.Range("K122").Value = 6410.099999999999
$ws.Range("M122").Value = -3960.099999999999
$ws.Range("H132").Value = 1539905.2
$ws.Range("I132").Value = 2748327.8
$ws.Range("K132").Value = 8244983.399999999
$ws.Range("M132").Value = -8242453.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6127
$ws.Range("I16").Value = 6127
$ws.Range("K16").Value = 6127
$ws.Range("M16").Value = -5957
$ws.Range("H22").Value = 2294.4
$ws.Range("I22").Value = 1249.5
$ws.Range("J22").Value = 2991
$ws.Range("K22").Value = 1249.5
$ws.Range("L22").Value = 2991
$ws.Range("M22").Value = -954.5
$ws.Range("N22").Value = -3581
$ws.Range("H27").Value = 2294.4
$ws.Range("I27").Value = 1249.5
$ws.Range("J27").Value = 2991
$ws.Range("K27").Value = 1249.5
$ws.Range("L27").Value = 2991
$ws.Range("M27").Value = -1142.5
$ws.Range("N27").Value = -3205
$ws.Range("H93").Value = 30304018
$ws.Range("I93").Value = 885.8889
$ws.Range("K93").Value = 885.8889
$ws.Range("M93").Value = 362.1111
$ws.Range("H122").Value = 7208.2085
$ws.Range("J122").Value = 6324.875
$ws.Range("L122").Value = 18974.625
$ws.Range("N122").Value = -23874.625
$ws.Range("H136").Value = 2809.7
$ws.Range("I136").Value = 2715.3333
$ws.Range("K136").Value = 8145.999899999999
$ws.Range("M136").Value = -5595.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 28995.322
$ws.Range("I122").Value = 52826.734
$ws.Range("K122").Value = 158480.202
$ws.Range("M122").Value = -156030.202
$ws.Range("H132").Value = 1465.4822
$ws.Range("I132").Value = 1067.0769
$ws.Range("J132").Value = 2379.4707
$ws.Range("K132").Value = 3201.2307
$ws.Range("L132").Value = 7138.4121
$ws.Range("M132").Value = -671.2307000000001
$ws.Range("N132").Value = -12198.4121
$ws.Range("H136").Value = 1159.8064
$ws.Range("I136").Value = 861.76
$ws.Range("J136").Value = 2401.6667
$ws.Range("K136").Value = 2585.28
$ws.Range("L136").Value = 7205.000100000001
$ws.Range("M136").Value = -35.27999999999975
$ws.Range("N136").Value = -12305.0001
